# Split the "compoundStmt = "{" statement "}" ." run on the code-generation
# grammar-rule slide into three runs, changing "statement" to "statements"
# in the middle run:
#
#   ' = "{" statement "}" .'
#       -> ' = '  +  '"{" statements '  +  '"}" .'

$p = $ppt.ActivePresentation

$oldFragment = ' = "{" statement "}" .'
$partA       = ' = '
$partB       = '"{" statements '
$partC       = '"}" .'

# Locate the shape (on any slide) whose text contains the grammar rule.
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.Contains($oldFragment)) {
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not locate the shape containing '$oldFragment'"
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text

# 1-based start position (PowerPoint TextRange indices are 1-based) of the
# fragment we are rewriting.
$start = $fullText.IndexOf($oldFragment) + 1

$startA = $start
$startB = $startA + $partA.Length
$startC = $startB + $partB.Length

# Apply the edits from right to left so earlier (lower) character offsets
# stay valid while later ones are being rewritten.
$rangeC = $tr.Characters($startC, $partC.Length)
$rangeC.Text = $partC

$rangeB = $tr.Characters($startB, $partB.Length)
$rangeB.Text = $partB

$rangeA = $tr.Characters($startA, $partA.Length)
$rangeA.Text = $partA
